$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.134287357330322
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.061272025108337
